# Winterhold-cards.xlsx edit script
# Applies the changes described in the commit:
#   "Update version to 0.1.0 Remove some cards from the excel sheet"
#
# Summary of data changes on the "Destruction" worksheet:
#   - Row 9  (Tome Stack):      energy-cost code column (F) was blank, now "x"
#   - Row 13 (Chain Lightning): description changed to a Shock AoE effect
#   - Row 17 (Sparks):          description changed to apply Jumpy Lightning
#   - Row 21:                   gains the "Simon Says" card data (moved up from row 28)
#   - Rows 25-28 (HoseDown, OilUp, Ground, Simon Says): card entries removed;
#     their I/J reference-table cells (Spellweave/Singe/Jumpy Lightning) are kept

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9: Tome Stack now has an energy "code" of x ---
$ws.Range("F9").Value = "x"

# --- Row 13: Chain Lightning gets a new description ---
$ws.Range("C13").Value = "Deal 8 Shock damage to ALL enemies."

# --- Row 17: Sparks gets a new description ---
$ws.Range("C17").Value = "Apply 1(2) Jumpy Lightning."

# --- Row 21: Simon Says card moves here from row 28 ---
$ws.Range("A21").Value = "Simon Says"
$ws.Range("B21").Value = "Power"
$ws.Range("C21").Value = "Each time you deal spell damage, this power chooses Fire, Frost, or Shock randomly. If the next damage you deal is of that type, gain 1 Strength. If not, lose all Strength."
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = "Rare"

# --- Rows 25-28: remove the HoseDown / OilUp / Ground / Simon Says(old)
#     card entries by clearing their contents (rather than deleting the
#     rows), which keeps the unrelated I/J "Combo system" reference rows
#     (Spellweave, Singe, Jumpy Lightning) intact and leaves other
#     formulas' cell ranges untouched ---
$ws.Range("A25:F28").ClearContents()

# --- Leave the selection where the last edit was made (matches the saved
#     workbook's cursor position) ---
$ws.Range("C14").Select() | Out-Null
